$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "21.809.67"
$ws.Range("E2").Value = "  -1.38%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.540.61"
$ws.Range("E3").Value = "  -0.98%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.64%  "

# Row 5 - USDC
$ws.Range("D5").Value = "1.005"

# Row 6 - BNB
$ws.Range("D6").Value = "289.50"
$ws.Range("E6").Value = "  +0.53%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.3942"
$ws.Range("E7").Value = "  +3.97%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3198"
$ws.Range("E8").Value = "  -2.78%  "

# Row 9 - OKB
$ws.Range("D9").Value = "43.24"
$ws.Range("E9").Value = "  +0.09%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.07180"
$ws.Range("E10").Value = "  -2.51%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  -6.73%  "

# Row 12 - BinanceUSD
$ws.Range("D12").Value = "1.006"
$ws.Range("E12").Value = "  +0.63%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "5.628"
$ws.Range("E13").Value = "  -3.34%  "

# Row 14 - Solana
$ws.Range("D14").Value = "18.55"
$ws.Range("E14").Value = "  -8.13%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "6.626"
$ws.Range("E15").Value = "  -3.00%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.545.18"
$ws.Range("E16").Value = "  -1.21%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.00001099"
$ws.Range("E17").Value = "  -0.40%  "

# Row 18 - TRON
$ws.Range("D18").Value = "0.06554"
$ws.Range("E18").Value = "  -0.94%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "83.28"
$ws.Range("E19").Value = "  -3.04%  "

# Row 20 - Dai
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.54%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.142"
$ws.Range("E21").Value = "  -3.95%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -5.17%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "11.01"
$ws.Range("E23").Value = "  -6.02%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.382"
$ws.Range("E24").Value = "  +2.89%  "

# Row 25 - WrappedBTC
$ws.Range("D25").Value = "21.841.97"

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "2.374"
$ws.Range("E26").Value = "  -6.06%  "

# Row 27 - Monero
$ws.Range("D27").Value = "145.15"
$ws.Range("E27").Value = "  -3.61%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "18.37"
$ws.Range("E28").Value = "  -3.98%  "

# Row 29 - HuobiToken
$ws.Range("D29").Value = "4.857"
$ws.Range("E29").Value = "  -1.12%  "

# Row 30 - WrappedliquidstakedEther2.0
$ws.Range("D30").Value = "1.717.01"
$ws.Range("E30").Value = "  -0.88%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "117.20"
$ws.Range("E31").Value = "  -3.72%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "0.9678"
$ws.Range("E32").Value = "  -10.55%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "5.877"
$ws.Range("E33").Value = "  -1.52%  "

# Row 34 - Stellar
$ws.Range("D34").Value = "0.08232"
$ws.Range("E34").Value = "  +0.07%  "

# Row 35 - FraxShare
$ws.Range("D35").Value = "8.965"
$ws.Range("E35").Value = "  -4.01%  "

# Row 36 - WEMIXTOKEN
$ws.Range("D36").Value = "1.556"
$ws.Range("E36").Value = "  -15.93%  "

# Row 37 - Hedera
$ws.Range("D37").Value = "0.06060"
$ws.Range("E37").Value = "  -3.19%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -4.34%  "

# Row 39 - InternetComputer(DFINITY)
$ws.Range("D39").Value = "5.109"
$ws.Range("E39").Value = "  -3.58%  "

# Row 40 & 41 - Algorand/TrustWalletToken swap places
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.196"
$ws.Range("E40").Value = "  -4.68%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.2030"
$ws.Range("E41").Value = "  -6.13%  "

# Row 42 - Frax
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  +0.52%  "

# Row 43 - Aptos
$ws.Range("D43").Value = "10.63"
$ws.Range("E43").Value = "  -3.91%  "

# Row 44 - TheSandbox
$ws.Range("D44").Value = "0.5761"
$ws.Range("E44").Value = "  -4.97%  "

# Row 45 - PancakeSwap
$ws.Range("D45").Value = "3.757"
$ws.Range("E45").Value = "  +0.50%  "

# Row 46 - EnergySwap
$ws.Range("E46").Value = "  -5.50%  "

# Row 47 - Decentraland
$ws.Range("D47").Value = "0.5542"
$ws.Range("E47").Value = "  -5.36%  "

# Row 48 - Quant
$ws.Range("D48").Value = "116.95"
$ws.Range("E48").Value = "  -4.48%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "1.860"
$ws.Range("E49").Value = "  -6.63%  "

# Row 50 - EOS
$ws.Range("D50").Value = "1.130"
$ws.Range("E50").Value = "  -4.01%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.06747"
$ws.Range("E51").Value = "  -3.91%  "
